# "Generate Report for Handback"
#
# For each locale sheet (zh-cn, de-de):
#   - Column F ("Latest Target File")   gets the same file + hyperlink as column A
#     ("Source File Name") for that row.
#   - Column G ("Latest Handback File") gets the same file + hyperlink as column D
#     ("Latest Handoff File") for that row.
#   - Column H ("Latest Handback DateTime") is stamped with the handback time.
# The overall Status text moves from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview + both
# locale sheets).

$wb = $excel.ActiveWorkbook

function Get-LinkAddress($ws, $cellAddress) {
    foreach ($link in $ws.Hyperlinks) {
        if ($link.Range.Address() -eq $cellAddress) {
            return $link.Address
        }
    }
    return $null
}

function Add-MatchingHyperlink($ws, $targetCell, $sourceCellAddress, $displayText) {
    $url = Get-LinkAddress $ws $sourceCellAddress
    if ($url -ne $null) {
        $ws.Hyperlinks.Add($ws.Range($targetCell), $url, "", "", $displayText) | Out-Null
    } else {
        $ws.Range($targetCell).Value = $displayText
    }
    # Match the look of the existing hyperlink cells (single underline, the
    # workbook's custom hyperlink blue).
    $ws.Range($targetCell).Font.Underline = 2
    $ws.Range($targetCell).Font.Color = 15570276
}

# --- 1. Flip the "Ready for handoff" status everywhere to the handback status ---
$statusCells = @(
    @{Sheet="Overview"; Cells=@("B2","C2","B3","C3")},
    @{Sheet="zh-cn";     Cells=@("C2","C3")},
    @{Sheet="de-de";     Cells=@("C2","C3")}
)
foreach ($entry in $statusCells) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($addr in $entry.Cells) {
        $ws.Range($addr).Value = "Handed back: in sync with en-US"
    }
}

# --- 2. Populate the handback columns (F/G/H) on each locale sheet ---
$handbackTimes = @{ "zh-cn" = "2016-03-17 10:15:10"; "de-de" = "2016-03-17 10:15:20" }

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($r in @(2, 3)) {
        $aAddr = "`$A`$$r"
        $dAddr = "`$D`$$r"

        $aDisplay = $ws.Range("A$r").Value2
        $dDisplay = $ws.Range("D$r").Value2

        Add-MatchingHyperlink $ws "F$r" $aAddr $aDisplay
        Add-MatchingHyperlink $ws "G$r" $dAddr $dDisplay

        $ws.Range("H$r").Value = $handbackTimes[$sheetName]
    }
}
